{"js": "// Add the new `Balance` column to the `Account` table, directly after\n// the `PersonId` column definition (and before `BankId`).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst personIdFieldPattern = /^\\s*`PersonId`\\s+int\\(/;\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (personIdFieldPattern.test(paragraphs.items[i].text)) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the `PersonId` column definition in the Account table.\");\n}\n\ntarget.insertParagraph(\"  `Balance` int(10) NOT NULL,\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add the new `Balance` column to the `Account` table, directly after\n# the `PersonId` column definition (and before `BankId`).\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like '*`PersonId`* int(*') {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate the ``PersonId`` column definition in the Account table.\"\n}\n\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"  ``Balance`` int(10) NOT NULL,\"\n"}
